# Commit: "Added  FolderManager::ComputeParentFolderSize, more test cases"
#
# 1. Insert a new blank worksheet "C.Users" before the first sheet
#    ("C.Users.acarz"), becoming the new first tab.
# 2. Append new test-run rows (19-25, 30-37) to the "C.Users.acarz" sheet.
# 3. Update the saved selections: "C.Users.Default" loses the active tab /
#    its selection moves to B26; "C.Users.acarz" becomes the active tab
#    with its selection on C28.

$wb = $excel.ActiveWorkbook

# --- 1. Insert the new "C.Users" sheet in front of everything else -------
$acarzSheet = $wb.Worksheets.Item("C.Users.acarz")
$newSheet = $wb.Worksheets.Add($acarzSheet)
$newSheet.Name = "C.Users"

# --- 2. Fill in the new rows on "C.Users.acarz" ---------------------------
$ws = $wb.Worksheets.Item("C.Users.acarz")

$ws.Range("C19").Value = 'Before: 2025-04-07 02:06:53'
$ws.Range("C20").Value = 'After: 2025-04-07 02:07:42 Folder size: 233462132864'
$ws.Range("C21").Value = 'Number of folders: 35'
$ws.Range("C22").Value = 'After2: 2025-04-07 02:07:42'
$ws.Range("C23").Value = 'Number of files: 13'
$ws.Range("C24").Value = 'After3: 2025-04-07 02:07:42'
$ws.Range("C25").Value = 'Folder C:\Users\acarz size: 233462132864 last checked: 1969-12-31 23:59:59 last modified: 1969-12-31 23:59:59'
$ws.Range("E25").Value = 'Folders::GetInstance().ComputeFolderSizeInternally("C:\\Users\\acarz");'
$ws.Range("E25").Font.Bold = $true
$ws.Range("E20").Value = 'Windows explorer shows 284,193,604,168'

$ws.Range("C31").Value = 'Before: 2025-04-07 02:35:54'
$ws.Range("C32").Value = 'After: 2025-04-07 02:37:07 Folder size: 275836550185'
$ws.Range("C33").Value = 'Number of folders: 35'
$ws.Range("C34").Value = 'After2: 2025-04-07 02:37:07'
$ws.Range("C35").Value = 'Number of files: 13'
$ws.Range("C36").Value = 'After3: 2025-04-07 02:37:07'
$ws.Range("C37").Value = 'Folder C:\Users\acarz size: 275836550185 last checked: 1969-12-31 23:59:59 last modified: 1969-12-31 23:59:59'
$ws.Range("C30").Value = 'executed as admin :'

# --- 3. Update view/selection state ---------------------------------------
# "C.Users.Default" keeps a B26 selection but is no longer the active tab.
$defaultSheet = $wb.Worksheets.Item("C.Users.Default")
$defaultSheet.Range("B26").Select()

# "C.Users.acarz" becomes the active tab with selection C28 (set last so it
# wins the "active sheet" state).
$ws.Range("C28").Select()
